# Burn Down Chart update for Fradet (Richard):
#  - Log actual effort hours worked against his tasks (rows 4-6, column H/I)
#  - Leave the workbook with the "How To Use" sheet active/selected,
#    matching where the author ended up when they saved.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Burn Down Chart")
$ws2 = $wb.Worksheets.Item("How To Use")

# Richard logged 2 hours on Day 2 (task M2-2, row 4)
$ws1.Range("H4").Value = 2
# Richard logged 1 hour on Day 2 (task M2-3, row 5)
$ws1.Range("H5").Value = 1
# Richard logged 1 hour on Day 3 (task M2-4, row 6)
$ws1.Range("I6").Value = 1

# Move the selection on the burn-down sheet, then switch the active tab
# over to "How To Use" before saving (matches the author's final view).
$ws1.Range("J7").Select() | Out-Null
$ws2.Activate() | Out-Null
